$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing fixtures: swapped/rotated rows (130/132, 134/136, 143/144/145) ---
# Row 130
$ws.Range('B130').Value = 7483189
$ws.Range('F130').Value = 'Independiente del Valle'
$ws.Range('G130').Value = 'Orense'
$ws.Range('H130').Value = 2
$ws.Range('I130').Value = 2
$ws.Range('J130').Value = 'D'
$ws.Range('K130').Value = 1.4
$ws.Range('L130').Value = 4.75
$ws.Range('M130').Value = 7
$ws.Range('N130').Value = 1.4
$ws.Range('O130').Value = 4.5
$ws.Range('P130').Value = 8
$ws.Range('Q130').Value = -1.25
$ws.Range('R130').Value = 1.875
$ws.Range('S130').Value = 1.925
$ws.Range('T130').Value = 2.5
$ws.Range('U130').Value = 1.925
$ws.Range('V130').Value = 1.875
$ws.Range('W130').Value = -1
$ws.Range('X130').Value = 3.5
$ws.Range('Y130').Value = -1
$ws.Range('Z130').Value = -1
$ws.Range('AA130').Value = 0.925
$ws.Range('AB130').Value = 0.925
$ws.Range('AC130').Value = -1

# Row 132
$ws.Range('B132').Value = 7483081
$ws.Range('F132').Value = 'Deportivo Cuenca'
$ws.Range('G132').Value = 'El Nacional'
$ws.Range('H132').Value = 1
$ws.Range('I132').Value = 0
$ws.Range('J132').Value = 'H'
$ws.Range('K132').Value = 2.75
$ws.Range('L132').Value = 3.25
$ws.Range('M132').Value = 2.55
$ws.Range('N132').Value = 3
$ws.Range('O132').Value = 3.3
$ws.Range('P132').Value = 2.3
$ws.Range('Q132').Value = 0.25
$ws.Range('R132').Value = 1.825
$ws.Range('S132').Value = 1.975
$ws.Range('T132').Value = 2.75
$ws.Range('U132').Value = 2
$ws.Range('V132').Value = 1.8
$ws.Range('W132').Value = 2
$ws.Range('X132').Value = -1
$ws.Range('Y132').Value = -1
$ws.Range('Z132').Value = 0.825
$ws.Range('AA132').Value = -1
$ws.Range('AB132').Value = -1
$ws.Range('AC132').Value = 0.8

# Row 134
$ws.Range('B134').Value = 7483188
$ws.Range('F134').Value = 'Gualaceo SC'
$ws.Range('G134').Value = 'Emelec'
$ws.Range('H134').Value = 0
$ws.Range('I134').Value = 2
$ws.Range('J134').Value = 'A'
$ws.Range('K134').Value = 3.6
$ws.Range('L134').Value = 3.3
$ws.Range('M134').Value = 2.05
$ws.Range('N134').Value = 2.6
$ws.Range('O134').Value = 3.25
$ws.Range('P134').Value = 2.75
$ws.Range('Q134').Value = 0
$ws.Range('R134').Value = 1.8
$ws.Range('S134').Value = 2
$ws.Range('T134').Value = 2.5
$ws.Range('U134').Value = 1.975
$ws.Range('V134').Value = 1.825
$ws.Range('W134').Value = -1
$ws.Range('X134').Value = -1
$ws.Range('Y134').Value = 1.75
$ws.Range('Z134').Value = -1
$ws.Range('AA134').Value = 1
$ws.Range('AB134').Value = -1
$ws.Range('AC134').Value = 0.825

# Row 136
$ws.Range('B136').Value = 7482867
$ws.Range('F136').Value = 'Cumbaya FC'
$ws.Range('G136').Value = 'LDU Quito'
$ws.Range('H136').Value = 1
$ws.Range('I136').Value = 2
$ws.Range('J136').Value = 'A'
$ws.Range('K136').Value = 5.25
$ws.Range('L136').Value = 3.75
$ws.Range('M136').Value = 1.65
$ws.Range('N136').Value = 9
$ws.Range('O136').Value = 4.5
$ws.Range('P136').Value = 1.363
$ws.Range('Q136').Value = 1.25
$ws.Range('R136').Value = 1.975
$ws.Range('S136').Value = 1.825
$ws.Range('T136').Value = 2.5
$ws.Range('U136').Value = 1.825
$ws.Range('V136').Value = 1.975
$ws.Range('W136').Value = -1
$ws.Range('X136').Value = -1
$ws.Range('Y136').Value = 0.363
$ws.Range('Z136').Value = 0.4875
$ws.Range('AA136').Value = -0.5
$ws.Range('AB136').Value = 0.825
$ws.Range('AC136').Value = -1

# Row 143
$ws.Range('B143').Value = 7528857
$ws.Range('F143').Value = 'Universidad Catolica del Ecuador'
$ws.Range('G143').Value = 'Barcelona Guayaquil'
$ws.Range('H143').Value = 0
$ws.Range('I143').Value = 1
$ws.Range('J143').Value = 'A'
$ws.Range('K143').Value = 1.533
$ws.Range('L143').Value = 4
$ws.Range('M143').Value = 5.5
$ws.Range('N143').Value = 1.5
$ws.Range('O143').Value = 4.333
$ws.Range('P143').Value = 5.25
$ws.Range('Q143').Value = -1
$ws.Range('R143').Value = 1.8
$ws.Range('S143').Value = 2
$ws.Range('T143').Value = 3
$ws.Range('U143').Value = 1.975
$ws.Range('V143').Value = 1.825
$ws.Range('W143').Value = -1
$ws.Range('X143').Value = -1
$ws.Range('Y143').Value = 4.25
$ws.Range('Z143').Value = -1
$ws.Range('AA143').Value = 1
$ws.Range('AB143').Value = -1
$ws.Range('AC143').Value = 0.825

# Row 144
$ws.Range('B144').Value = 7528858
$ws.Range('F144').Value = 'Orense'
$ws.Range('G144').Value = 'SD Aucas'
$ws.Range('H144').Value = 1
$ws.Range('I144').Value = 2
$ws.Range('J144').Value = 'A'
$ws.Range('K144').Value = 2.2
$ws.Range('L144').Value = 3.2
$ws.Range('M144').Value = 3.2
$ws.Range('N144').Value = 1.95
$ws.Range('O144').Value = 3.2
$ws.Range('P144').Value = 3.8
$ws.Range('Q144').Value = -0.5
$ws.Range('R144').Value = 1.95
$ws.Range('S144').Value = 1.85
$ws.Range('T144').Value = 2.25
$ws.Range('U144').Value = 1.85
$ws.Range('V144').Value = 1.95
$ws.Range('W144').Value = -1
$ws.Range('X144').Value = -1
$ws.Range('Y144').Value = 2.8
$ws.Range('Z144').Value = -1
$ws.Range('AA144').Value = 0.8500000000000001
$ws.Range('AB144').Value = 0.8500000000000001
$ws.Range('AC144').Value = -1

# Row 145
$ws.Range('B145').Value = 7528852
$ws.Range('F145').Value = 'Delfin SC'
$ws.Range('G145').Value = 'Tecnico Universitario'
$ws.Range('H145').Value = 2
$ws.Range('I145').Value = 2
$ws.Range('J145').Value = 'D'
$ws.Range('K145').Value = 2.1
$ws.Range('L145').Value = 3.4
$ws.Range('M145').Value = 3.1
$ws.Range('N145').Value = 2.1
$ws.Range('O145').Value = 3.4
$ws.Range('P145').Value = 3.1
$ws.Range('Q145').Value = -0.25
$ws.Range('R145').Value = 1.8
$ws.Range('S145').Value = 2
$ws.Range('T145').Value = 2.25
$ws.Range('U145').Value = 1.9
$ws.Range('V145').Value = 1.9
$ws.Range('W145').Value = -1
$ws.Range('X145').Value = 2.4
$ws.Range('Y145').Value = -1
$ws.Range('Z145').Value = -0.5
$ws.Range('AA145').Value = 0.5
$ws.Range('AB145').Value = 0.8999999999999999
$ws.Range('AC145').Value = -1

# --- Fill in results for previously-unplayed fixtures (rows 157-162) + rewrite row 163 ---
# Row 157
$ws.Range('H157').Value = 1
$ws.Range('I157').Value = 0
$ws.Range('J157').Value = 'H'
$ws.Range('K157').Value = 2.6
$ws.Range('L157').Value = 3.2
$ws.Range('M157').Value = 2.5
$ws.Range('N157').Value = 1.333
$ws.Range('O157').Value = 4.75
$ws.Range('P157').Value = 8
$ws.Range('Q157').Value = -1.5
$ws.Range('R157').Value = 2
$ws.Range('S157').Value = 1.8
$ws.Range('T157').Value = 2.75
$ws.Range('U157').Value = 1.95
$ws.Range('V157').Value = 1.85
$ws.Range('W157').Value = 0.333
$ws.Range('X157').Value = -1
$ws.Range('Y157').Value = -1
$ws.Range('Z157').Value = -1
$ws.Range('AA157').Value = 0.8
$ws.Range('AB157').Value = -1
$ws.Range('AC157').Value = 0.8500000000000001

# Row 158
$ws.Range('H158').Value = 3
$ws.Range('I158').Value = 0
$ws.Range('J158').Value = 'H'
$ws.Range('K158').Value = 2.5
$ws.Range('L158').Value = 3.2
$ws.Range('M158').Value = 2.6
$ws.Range('N158').Value = 2.5
$ws.Range('O158').Value = 3.2
$ws.Range('P158').Value = 2.625
$ws.Range('Q158').Value = 0
$ws.Range('R158').Value = 1.85
$ws.Range('S158').Value = 1.95
$ws.Range('T158').Value = 2.25
$ws.Range('U158').Value = 1.85
$ws.Range('V158').Value = 1.95
$ws.Range('W158').Value = 1.5
$ws.Range('X158').Value = -1
$ws.Range('Y158').Value = -1
$ws.Range('Z158').Value = 0.8500000000000001
$ws.Range('AA158').Value = -1
$ws.Range('AB158').Value = 0.8500000000000001
$ws.Range('AC158').Value = -1

# Row 159
$ws.Range('H159').Value = 3
$ws.Range('I159').Value = 1
$ws.Range('J159').Value = 'H'
$ws.Range('K159').Value = 1.666
$ws.Range('L159').Value = 3.75
$ws.Range('M159').Value = 4.5
$ws.Range('N159').Value = 1.7
$ws.Range('O159').Value = 3.6
$ws.Range('P159').Value = 4.5
$ws.Range('Q159').Value = -0.75
$ws.Range('R159').Value = 1.95
$ws.Range('S159').Value = 1.85
$ws.Range('T159').Value = 2.25
$ws.Range('U159').Value = 1.8
$ws.Range('V159').Value = 2
$ws.Range('W159').Value = 0.7
$ws.Range('X159').Value = -1
$ws.Range('Y159').Value = -1
$ws.Range('Z159').Value = 0.95
$ws.Range('AA159').Value = -1
$ws.Range('AB159').Value = 0.8
$ws.Range('AC159').Value = -1

# Row 160
$ws.Range('H160').Value = 3
$ws.Range('I160').Value = 1
$ws.Range('J160').Value = 'H'
$ws.Range('K160').Value = 2
$ws.Range('L160').Value = 3.4
$ws.Range('M160').Value = 3.4
$ws.Range('N160').Value = 1.95
$ws.Range('O160').Value = 3.4
$ws.Range('P160').Value = 3.6
$ws.Range('Q160').Value = -0.5
$ws.Range('R160').Value = 1.975
$ws.Range('S160').Value = 1.825
$ws.Range('T160').Value = 2.25
$ws.Range('U160').Value = 1.9
$ws.Range('V160').Value = 1.9
$ws.Range('W160').Value = 0.95
$ws.Range('X160').Value = -1
$ws.Range('Y160').Value = -1
$ws.Range('Z160').Value = 0.9750000000000001
$ws.Range('AA160').Value = -1
$ws.Range('AB160').Value = 0.8999999999999999
$ws.Range('AC160').Value = -1

# Row 161
$ws.Range('H161').Value = 1
$ws.Range('I161').Value = 1
$ws.Range('J161').Value = 'D'
$ws.Range('K161').Value = 1.833
$ws.Range('L161').Value = 3.25
$ws.Range('M161').Value = 4.2
$ws.Range('N161').Value = 1.75
$ws.Range('O161').Value = 3.3
$ws.Range('P161').Value = 4.75
$ws.Range('Q161').Value = -0.75
$ws.Range('R161').Value = 2
$ws.Range('S161').Value = 1.8
$ws.Range('T161').Value = 2.25
$ws.Range('U161').Value = 2
$ws.Range('V161').Value = 1.8
$ws.Range('W161').Value = -1
$ws.Range('X161').Value = 2.3
$ws.Range('Y161').Value = -1
$ws.Range('Z161').Value = -1
$ws.Range('AA161').Value = 0.8
$ws.Range('AB161').Value = -0.5
$ws.Range('AC161').Value = 0.4

# Row 162
$ws.Range('H162').Value = 1
$ws.Range('I162').Value = 3
$ws.Range('J162').Value = 'A'
$ws.Range('K162').Value = 5
$ws.Range('L162').Value = 3.6
$ws.Range('M162').Value = 1.615
$ws.Range('N162').Value = 8.5
$ws.Range('O162').Value = 4.2
$ws.Range('P162').Value = 1.363
$ws.Range('Q162').Value = 1.25
$ws.Range('R162').Value = 1.9
$ws.Range('S162').Value = 1.9
$ws.Range('T162').Value = 2.25
$ws.Range('U162').Value = 1.775
$ws.Range('V162').Value = 2.025
$ws.Range('W162').Value = -1
$ws.Range('X162').Value = -1
$ws.Range('Y162').Value = 0.363
$ws.Range('Z162').Value = -1
$ws.Range('AA162').Value = 0.8999999999999999
$ws.Range('AB162').Value = 0.7749999999999999
$ws.Range('AC162').Value = -1

# Row 163
$ws.Range('B163').Value = 7773062
$ws.Range('E163').Value = 45366.875
$ws.Range('F163').Value = 'Independiente del Valle'
$ws.Range('G163').Value = 'Cumbaya FC'
$ws.Range('K163').Value = 1.2
$ws.Range('L163').Value = 6
$ws.Range('M163').Value = 13
$ws.Range('N163').Value = 1.222
$ws.Range('O163').Value = 6
$ws.Range('P163').Value = 12
$ws.Range('Q163').Value = -1.75
$ws.Range('R163').Value = 1.925
$ws.Range('S163').Value = 1.875
$ws.Range('T163').Value = 2.75
$ws.Range('U163').Value = 1.9
$ws.Range('V163').Value = 1.9
$ws.Range('W163').Value = 0
$ws.Range('X163').Value = 0
$ws.Range('Y163').Value = 0
$ws.Range('Z163').Value = 0
$ws.Range('AA163').Value = 0

# --- Append new fixtures (rows 164-167) ---
# Row 164
$ws.Range('A163').Copy()
$ws.Range('A164').PasteSpecial(-4122)
$ws.Range('E163').Copy()
$ws.Range('E164').PasteSpecial(-4122)
$ws.Range('A164').Value = 162
$ws.Range('B164').Value = 7773471
$ws.Range('C164').Value = 'Ecuador LigaPro Serie A'
$ws.Range('D164').Value = 'Ecuador LigaPro Serie A'
$ws.Range('E164').Value = 45367.72916666666
$ws.Range('F164').Value = 'LDU Quito'
$ws.Range('G164').Value = 'SD Aucas'
$ws.Range('K164').Value = 1.533
$ws.Range('L164').Value = 4
$ws.Range('M164').Value = 6
$ws.Range('N164').Value = 1.533
$ws.Range('O164').Value = 4
$ws.Range('P164').Value = 6
$ws.Range('Q164').Value = -1
$ws.Range('R164').Value = 1.9
$ws.Range('S164').Value = 1.9
$ws.Range('T164').Value = 2.75
$ws.Range('U164').Value = 1.925
$ws.Range('V164').Value = 1.875
$ws.Range('W164').Value = 0
$ws.Range('X164').Value = 0
$ws.Range('Y164').Value = 0
$ws.Range('Z164').Value = 0
$ws.Range('AA164').Value = 0

# Row 165
$ws.Range('A163').Copy()
$ws.Range('A165').PasteSpecial(-4122)
$ws.Range('E163').Copy()
$ws.Range('E165').PasteSpecial(-4122)
$ws.Range('A165').Value = 163
$ws.Range('B165').Value = 7773473
$ws.Range('C165').Value = 'Ecuador LigaPro Serie A'
$ws.Range('D165').Value = 'Ecuador LigaPro Serie A'
$ws.Range('E165').Value = 45367.83333333334
$ws.Range('F165').Value = 'Barcelona Guayaquil'
$ws.Range('G165').Value = 'Orense'
$ws.Range('K165').Value = 1.615
$ws.Range('L165').Value = 3.6
$ws.Range('M165').Value = 6
$ws.Range('N165').Value = 1.4
$ws.Range('O165').Value = 4.333
$ws.Range('P165').Value = 9
$ws.Range('Q165').Value = -1.25
$ws.Range('R165').Value = 1.85
$ws.Range('S165').Value = 1.95
$ws.Range('T165').Value = 2.5
$ws.Range('U165').Value = 1.85
$ws.Range('V165').Value = 1.95
$ws.Range('W165').Value = 0
$ws.Range('X165').Value = 0
$ws.Range('Y165').Value = 0
$ws.Range('Z165').Value = 0
$ws.Range('AA165').Value = 0

# Row 166
$ws.Range('A163').Copy()
$ws.Range('A166').PasteSpecial(-4122)
$ws.Range('E163').Copy()
$ws.Range('E166').PasteSpecial(-4122)
$ws.Range('A166').Value = 164
$ws.Range('B166').Value = 7773472
$ws.Range('C166').Value = 'Ecuador LigaPro Serie A'
$ws.Range('D166').Value = 'Ecuador LigaPro Serie A'
$ws.Range('E166').Value = 45368.72916666666
$ws.Range('F166').Value = 'Delfin SC'
$ws.Range('G166').Value = 'Tecnico Universitario'
$ws.Range('K166').Value = 2.05
$ws.Range('L166').Value = 3.2
$ws.Range('M166').Value = 3.75
$ws.Range('N166').Value = 2.2
$ws.Range('O166').Value = 3.2
$ws.Range('P166').Value = 3.4
$ws.Range('Q166').Value = -0.25
$ws.Range('R166').Value = 1.875
$ws.Range('S166').Value = 1.925
$ws.Range('T166').Value = 2.25
$ws.Range('U166').Value = 2.025
$ws.Range('V166').Value = 1.775
$ws.Range('W166').Value = 0
$ws.Range('X166').Value = 0
$ws.Range('Y166').Value = 0
$ws.Range('Z166').Value = 0
$ws.Range('AA166').Value = 0

# Row 167
$ws.Range('A163').Copy()
$ws.Range('A167').PasteSpecial(-4122)
$ws.Range('E163').Copy()
$ws.Range('E167').PasteSpecial(-4122)
$ws.Range('A167').Value = 165
$ws.Range('B167').Value = 7773470
$ws.Range('C167').Value = 'Ecuador LigaPro Serie A'
$ws.Range('D167').Value = 'Ecuador LigaPro Serie A'
$ws.Range('E167').Value = 45368.83333333334
$ws.Range('F167').Value = 'El Nacional'
$ws.Range('G167').Value = 'Emelec'
$ws.Range('K167').Value = 2.2
$ws.Range('L167').Value = 3.3
$ws.Range('M167').Value = 3.2
$ws.Range('N167').Value = 2.2
$ws.Range('O167').Value = 3.3
$ws.Range('P167').Value = 3.2
$ws.Range('Q167').Value = -0.25
$ws.Range('R167').Value = 1.925
$ws.Range('S167').Value = 1.875
$ws.Range('T167').Value = 2.5
$ws.Range('U167').Value = 1.95
$ws.Range('V167').Value = 1.85
$ws.Range('W167').Value = 0
$ws.Range('X167').Value = 0
$ws.Range('Y167').Value = 0
$ws.Range('Z167').Value = 0
$ws.Range('AA167').Value = 0

$excel.CutCopyMode = $false
